$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

$ws.Range("C2").Value = 16
$ws.Range("C3").Value = 7
$ws.Range("C5").Value = 20
$ws.Range("C6").Value = 14
$ws.Range("C7").Value = 21
$ws.Range("C8").Value = 19
$ws.Range("C10").Value = 21
$ws.Range("C11").Value = 13
$ws.Range("C12").Value = 12
$ws.Range("C15").Value = 14
$ws.Range("C16").Value = 15
$ws.Range("C17").Value = 18
$ws.Range("B18").Value = "<unin>"
$ws.Range("C18").Value = 16
